$wb = $excel.ActiveWorkbook

# Helper: write literal text (not auto-converted to a date/number) into a
# cell by staging it through a scratch cell that is explicitly formatted as
# Text, then pasting only the *values* (not formats) into the destination.
function Set-LiteralText($sheet, $destAddress, $text) {
    $scratch = $sheet.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $sheet.Range($destAddress).PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
    $excel.CutCopyMode = $false
}

$newHeaderDate = "2020-04-23"   # new forecast-date column (V)
$newRowDate    = "2020-05-07"   # new observation row (34)

foreach ($sheetName in @("cases", "deaths")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- New column V: header in row 1 ---------------------------------
    Set-LiteralText $ws "V1" $newHeaderDate

    # Materialize the blank placeholder cells V2:V20 (present-but-empty,
    # matching the existing sparse-triangle layout). Touching a benign
    # formatting no-op forces the engine to keep the cell in the sheet.
    $ws.Range("V2:V20").Font.Bold = $false

    # --- New column V: numeric forecast values for rows 21-34 ----------
    $colV = 22   # column V

    # --- Row 20 gains an "Observed" value in column B -------------------
    if ($sheetName -eq "cases") {
        $ws.Cells.Item(20, 2).Value = 49492
    } else {
        $ws.Cells.Item(20, 2).Value = 3313
    }

    if ($sheetName -eq "cases") {
        $vVals = @(52655,55559,58540,62691,66235,70468,76254,81056,87305,91304,98977,104427,109459,114432)
    } else {
        $vVals = @(3500,3728,3964,4297,4584,4930,5410,5813,6343,6686,7352,7830,8275,8719)
    }
    $r = 21
    foreach ($v in $vVals) {
        $ws.Cells.Item($r, $colV).Value = $v
        $r = $r + 1
    }

    # --- New row 34: observation date + forecast value ------------------
    Set-LiteralText $ws "A34" $newRowDate
    $ws.Range("B34:U34").Font.Bold = $false
    # V34 already set above as the last element of $vVals (row 34)
}
